function Remove-ParaByText {
    param($d, $text)
    $r = $d.Content
    $found = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $text"
        return
    }
    $r.Expand(4)
    $r.Delete()
}

$d = $word.ActiveDocument

# Remove obsolete / completed "Do zrobienia" and "Pomysły" items, plus the
# whole "Sugestie" section (heading + its two bullet items).
Remove-ParaByText $d "Zaktualizować temat rekurencji"
Remove-ParaByText $d "Utworzyć wersję prezentacji z małpą"
Remove-ParaByText $d "Dodać rozwiązanie zadania ewaluacyjnego"
Remove-ParaByText $d "Rekurencja – dodać przykład"
Remove-ParaByText $d "Zadania dla chętnych"
Remove-ParaByText $d "Sugestie"
Remove-ParaByText $d "Wersje plików w formacie"

# Mark "Dodać prezentację ..." as the new last active item in "Do zrobienia":
# append a _GoBack bookmark right after its text (mirrors where Word leaves
# its last-edit marker).
$r3 = $d.Content
$r3.Find.Execute("Dodać prezentację", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r3.Expand(4)
$r3.MoveEnd(1, -1)
$r3.Collapse(0)
$r3.InsertAfter("~")
$markRange = $d.Content
$markRange.Find.Execute("~", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $markRange)
$markRange.Text = ""

# The former "Wykorzystanie także języka Python ..." paragraph (last item of
# the now-removed "Sugestie" section) becomes an empty, unnumbered paragraph.
$r4 = $d.Content
$r4.Find.Execute("Wykorzystanie także języka", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r4.Expand(4)
$r4.Delete()
$emptyParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Akapitzlist"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r4.InsertXML($emptyParaXml)
